# Word COM-interop script implementing the commit:
# "setup font awesome install via grunt and doc revisions"
#
# This script performs a sequence of targeted Find & Replace operations
# against $word.ActiveDocument to reproduce the author's prose revisions
# throughout the paper, plus a couple of structural tweaks (the _GoBack
# bookmark relocation and the cached PAGE field text in the footer).

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $result = $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $result) {
        Write-Host "WARNING: find failed for: $find"
    } else {
        Write-Host "OK: $find"
    }
}

# 1. "deepen" -> "broaden"; "technology" -> "block cipher"
Replace-Text "I decided to take on the challenge of writing the AES algorithm in order to deepen" "I decided to take on the challenge of writing the AES algorithm in order to broaden"
Replace-Text "my understanding of the technology," "my understanding of the block cipher,"

# 2. "16-byte lengths" -> "16-byte blocks"
Replace-Text "cided to standardize on 16-byte lengths" "cided to standardize on 16-byte blocks"

# 3. "many requirements... elaborated in deta" -> "numerous requirements... elaborated in much greater deta"
Replace-Text "understanding of many requirements, which are elaborated in deta" "understanding of numerous requirements, which are elaborated in much greater deta"

# 4. "rounds process performs the encryption" -> "rounds perform the encryption"
Replace-Text "process performs the encryption" "perform the encryption"

# 5. "I wrote the application in Javascript" -> "I wrote the AES algorithm in Javascript"
Replace-Text "I wrote the application in Javascript" "I wrote the AES algorithm in Javascript"

# 6. Drop "Apple, " from the companies list
Replace-Text "highly optimized by companies like Apple, Google and Microsoft." "highly optimized by companies like Google and Microsoft."

# 7. Drop "though the decimal values were "
Replace-Text "this implementation (though the decimal values were actually stored as 32-bit values)." "this implementation (actually stored as 32-bit values)."

# 8. "implementation process" -> "development process" (first challenge intro)
Replace-Text "I faced two challenges during the implementation process." "I faced two challenges during the development process."

# 9. "simple, yet lacking clarity: any number multiplied" -> "deceptively simple: Any number multiplied"
Replace-Text "simple, yet lacking clarity: any number multiplied" "deceptively simple: Any number multiplied"

# 10. Comma + "a few tense hours" (this also removes the _GoBack bookmark that
#     previously sat inside this span; it gets re-added near the end below)
Replace-Text "despite this improper implementation whereas the test decryption vectors failed. This led to several hours" "despite this improper implementation, whereas the test decryption vectors failed. This led to a few tense hours"

# 11. "determine the defect" -> "identify the defect"
Replace-Text "Eventually, I was able to determine the defect" "Eventually, I was able to identify the defect"

# 12. Drop "really "
Replace-Text "The installation process is really easy and should take less than 10 minutes" "The installation process is easy and should take less than 10 minutes"

# 13. Drop " Funtastic" (application name simplification)
Replace-Text "Next, clone the AES Funtastic application by typing" "Next, clone the AES application by typing"

# 14. "interesting aspects ... algorithm is ... interactions" -> "interesting features ... application is ... rounds"
Replace-Text "One of the interesting aspects of the AES algorithm is tracking the transformation of the key and plaintext values as the algorithm progresses through multiple cipher interactions." "One of the interesting features of the AES application is tracking the transformation of the key and plaintext values as the algorithm progresses through multiple cipher rounds."

# 15. "unlike" -> "very different from"
Replace-Text "that was unlike anything" "that was very different from anything"

# 16. Drop the "Initially, ... build up the momentum required to complete it. " sentence
Replace-Text "function. Initially, the project looked overwhelming, but by breaking it into small parts and achieving small early successes, I was able to build up the momentum required to complete it. I look forward" "function. I look forward"

# Re-insert the _GoBack bookmark at its new location: right after
# "...for each individual cipher function. " and before "I look forward..."
# (the original _GoBack bookmark, which used to sit inside the text replaced
# in step 10 above, was implicitly dropped by that edit).
$rng = $d.Content
$found = $rng.Find.Execute("I look forward to potentially expanding the project", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found -and -not $d.Bookmarks.Exists("_GoBack")) {
    $bmRange = $d.Range($rng.Start, $rng.Start)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}

# Footer page-number field: cached display text "1" -> "8"
$footer = $d.Sections.Item(1).Footers.Item(1)
$footerResult = $footer.Range.Find.Execute("1", $true, $false, $false, $false, $false, $true, 1, $false, "8", 2)
Write-Host "Footer replace result: $footerResult"

Write-Host "All replacements attempted."
